$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 1-3: simple text replacements (100 -> 0M, 0 -> 0M, 11 -> 0M)
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# Insert 10 new rows right after row 3 (i.e. before the row that currently
# holds "0" at index 4), carrying the new per-iteration stats.
$newValues = @("13", "0.00003", "0.00004", "0.00003", "0.00001", "0.00003", "0.00003", "0.00004", "0.00042", "100.0")

$refRow = $t.Rows.Item(4)
for ($i = $newValues.Count - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
}

# The two tab-separated summary rows near the end collapse down to just
# their first value, and the trailing blank row picks up the last value.
$rowCount = $t.Rows.Count
$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = "100"
$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = "0"
$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = "11"
